$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply "Good" cell style to row 7 (A7:C7) to match rows 5, 6, 13
$ws.Range("A7:B7").Style = "Good"
$ws.Range("C7").Style = "Good"
$ws.Range("C7").HorizontalAlignment = -4108

# Apply "Good" cell style to row 10 (A10:C10)
$ws.Range("A10:B10").Style = "Good"
$ws.Range("C10").Style = "Good"
$ws.Range("C10").HorizontalAlignment = -4108

# Add the new task row 16 (column C keeps the default centered "style 2" already
# set on the whole column, same as the other non-highlighted rows)
$ws.Range("A16").Value = "Rūšiavimas skelbimų pagal radimo datą"
$ws.Range("B16").Value = "Ignas"
$ws.Range("C16").Value = "vidutinis"

# Update the view: scroll so row 4 is the top-left visible row, and select C8
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C8").Select()
